$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Target values per row for columns D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion), S (Precio $/Kg), T (Kg / unidad).
# This reflects a reshuffle of the weekly price rows 3..11.

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

$rows = @{
    3  = @{ D = $epoch.AddDays(44391); L = "Primera"; M = 15; N = 1500;  O = 1500;  P = 1500;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1500; T = 1 }
    4  = @{ D = $epoch.AddDays(44391); L = "Segunda"; M = 20; N = 1000;  O = 1000;  P = 1000;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1000; T = 1 }
    5  = @{ D = $epoch.AddDays(44309); L = "Primera"; M = 10; N = 1600;  O = 1600;  P = 1600;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1600; T = 1 }
    6  = @{ D = $epoch.AddDays(44400); L = "Primera"; M = 25; N = 1500;  O = 1500;  P = 1500;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1500; T = 1 }
    7  = @{ D = $epoch.AddDays(44343); L = "Primera"; M = 20; N = 1700;  O = 1700;  P = 1700;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1700; T = 1 }
    8  = @{ D = $epoch.AddDays(44371); L = "Primera"; M = 20; N = 1800;  O = 1800;  P = 1800;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1800; T = 1 }
    9  = @{ D = $epoch.AddDays(44371); L = "Segunda"; M = 30; N = 1200;  O = 1200;  P = 1200;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1200; T = 1 }
    10 = @{ D = $epoch.AddDays(44336); L = "Primera"; M = 10; N = 1500;  O = 1500;  P = 1500;  Q = "`$/kilo (en caja de 15 kilos)"; S = 1500; T = 1 }
    11 = @{ D = $epoch.AddDays(44195); L = "Primera"; M = 20; N = 15000; O = 15000; P = 15000; Q = "`$/bandeja 10 kilos";            S = 1500; T = 10 }
}

foreach ($r in $rows.Keys) {
    $v = $rows[$r]
    $ws.Cells.Item($r, 4).Value = $v.D
    $ws.Cells.Item($r, 12).Value = $v.L
    $ws.Cells.Item($r, 13).Value = $v.M
    $ws.Cells.Item($r, 14).Value = $v.N
    $ws.Cells.Item($r, 15).Value = $v.O
    $ws.Cells.Item($r, 16).Value = $v.P
    $ws.Cells.Item($r, 17).Value = $v.Q
    $ws.Cells.Item($r, 19).Value = $v.S
    $ws.Cells.Item($r, 20).Value = $v.T
}
